$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 137  # was 136
$ws.Range("F3").Value = 131  # was 132
$ws.Range("F4").Value = 910  # was 909
$ws.Range("F5").Value = 1078  # was 1077
$ws.Range("F6").Value = 1551  # was 1549
$ws.Range("F7").Value = 334  # was 333
$ws.Range("F8").Value = 678  # was 677
$ws.Range("F9").Value = 12259  # was 12224
$ws.Range("F11").Value = 2182  # was 2177
$ws.Range("F12").Value = 906  # was 904
$ws.Range("F13").Value = 260  # was 257
$ws.Range("F16").Value = 219  # was 214
$ws.Range("F17").Value = 277  # was 275
$ws.Range("F18").Value = 782  # was 780
$ws.Range("F20").Value = 303  # was 301
$ws.Range("F21").Value = 2928  # was 2924
$ws.Range("F22").Value = 763  # was 761
$ws.Range("F23").Value = 4186  # was 4152
$ws.Range("F24").Value = 4186  # was 4152
$ws.Range("F25").Value = 1129  # was 1123
$ws.Range("F26").Value = 866  # was 864
$ws.Range("F28").Value = 29  # was 27
$ws.Range("F30").Value = 1059  # was 1057
$ws.Range("F31").Value = 55  # was 53
$ws.Range("F32").Value = 105  # was 103
$ws.Range("F33").Value = 271  # was 270
$ws.Range("F36").Value = 31  # was 29
$ws.Range("F38").Value = 4438  # was 4430
$ws.Range("F39").Value = 17  # was 16
$ws.Range("F40").Value = 4562  # was 4554
$ws.Range("F41").Value = 5553  # was 5545
$ws.Range("F43").Value = 130  # was 129
$ws.Range("F44").Value = 77  # was 70
$ws.Range("F45").Value = 174  # was 173
$ws.Range("F46").Value = 332  # was 328
$ws.Range("F47").Value = 82  # was 81
$ws.Range("F50").Value = 133  # was 132

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4172  # was 4173
$ws.Range("F5").Value = 63  # was 61
$ws.Range("G5").Value = 108  # was 88
$ws.Range("F11").Value = 111  # was 109
$ws.Range("F13").Value = 1039  # was 1038

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 762  # was 761
$ws.Range("F3").Value = 464  # was 461
$ws.Range("F4").Value = 82  # was 81
$ws.Range("F5").Value = 10  # was 7

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 762  # was 761
$ws.Range("F3").Value = 464  # was 461
$ws.Range("F4").Value = 82  # was 81
$ws.Range("F5").Value = 137  # was 136
$ws.Range("F6").Value = 910  # was 909
$ws.Range("F7").Value = 1551  # was 1549
$ws.Range("F8").Value = 334  # was 333
$ws.Range("F9").Value = 678  # was 677
$ws.Range("F10").Value = 12259  # was 12224
$ws.Range("F11").Value = 2182  # was 2177
$ws.Range("F12").Value = 260  # was 257
$ws.Range("F14").Value = 277  # was 275
$ws.Range("F16").Value = 303  # was 301
$ws.Range("F17").Value = 2928  # was 2924
$ws.Range("F18").Value = 763  # was 761
$ws.Range("F19").Value = 63  # was 61
$ws.Range("G19").Value = 108  # was 88
$ws.Range("F20").Value = 4186  # was 4152
$ws.Range("F21").Value = 1129  # was 1123
$ws.Range("F23").Value = 866  # was 864
$ws.Range("F27").Value = 1059  # was 1057
$ws.Range("F28").Value = 55  # was 53
$ws.Range("F29").Value = 105  # was 103
$ws.Range("F30").Value = 271  # was 270
$ws.Range("F32").Value = 31  # was 29
$ws.Range("F33").Value = 4438  # was 4430
$ws.Range("F34").Value = 4562  # was 4554
$ws.Range("F37").Value = 130  # was 129
$ws.Range("F38").Value = 174  # was 173
$ws.Range("F39").Value = 332  # was 328
$ws.Range("F42").Value = 82  # was 81
$ws.Range("F49").Value = 133  # was 132

$wb.Save()